$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.1266753333333333
$ws.Range("H2").Value = 0.380026
$ws.Range("I2").Value = 0.003969958931547584
$ws.Range("J2").Value = 0.004041902127696777
$ws.Range("M2").Value = 121.928739
$ws.Range("N2").Value = 365.786217
$ws.Range("O2").Value = 0.2282232151508951
$ws.Range("P2").Value = 0.2419720431319445
$ws.Range("Q2").Value = 15.445363655738
$ws.Range("R2").Value = 139.008272901642
$ws.Range("S2").Value = 0.0009060367913748018
$ws.Range("T2").Value = 0.0009780273159781428
$ws.Range("G3").Value = 0.1266753333333333
$ws.Range("H3").Value = 0.380026
$ws.Range("I3").Value = 0.003969958931547584
$ws.Range("J3").Value = 0.004041902127696777
$ws.Range("O3").Value = 0.2768624053389947
$ws.Range("P3").Value = 0.2935413991166814
$ws.Range("Q3").Value = 18.73709705752666
$ws.Range("R3").Value = 168.63387351774
$ws.Range("S3").Value = 0.00109913237888529
$ws.Range("T3").Value = 0.001186465605656803
$ws.Range("G4").Value = 0.1266753333333333
$ws.Range("H4").Value = 0.380026
$ws.Range("I4").Value = 0.003969958931547584
$ws.Range("J4").Value = 0.004041902127696777
$ws.Range("M4").Value = 83.50496933333334
$ws.Range("N4").Value = 250.514908
$ws.Range("O4").Value = 0.1563025480180701
$ws.Range("P4").Value = 0.1657186665504434
$ws.Range("Q4").Value = 10.57801982528978
$ws.Range("R4").Value = 95.20217842760799
$ws.Range("S4").Value = 0.0006205146965279823
$ws.Range("T4").Value = 0.0006698186309293097
$ws.Range("G5").Value = 0.1266753333333333
$ws.Range("H5").Value = 0.380026
$ws.Range("I5").Value = 0.003969958931547584
$ws.Range("J5").Value = 0.004041902127696777
$ws.Range("M5").Value = 91.06846250000001
$ws.Range("N5").Value = 182.136925
$ws.Range("O5").Value = 0.1704597085236707
$ws.Range("P5").Value = 0.1204857969594293
$ws.Range("Q5").Value = 11.53612784334167
$ws.Range("R5").Value = 69.21676706005
$ws.Range("S5").Value = 0.0006767180423225444
$ws.Range("T5").Value = 0.0004869917990875592
$ws.Range("G6").Value = 0.1266753333333333
$ws.Range("H6").Value = 0.380026
$ws.Range("I6").Value = 0.003969958931547584
$ws.Range("J6").Value = 0.004041902127696777
$ws.Range("M6").Value = 89.83562999999999
$ws.Range("N6").Value = 269.50689
$ws.Range("O6").Value = 0.1681521229683693
$ws.Range("P6").Value = 0.1782820942415013
$ws.Range("Q6").Value = 11.37995837546
$ws.Range("R6").Value = 102.41962537914
$ws.Range("S6").Value = 0.0006675570224369653
$ws.Range("T6").Value = 0.0007205987760449614
$ws.Range("I7").Value = 0.9391988012463586
$ws.Range("J7").Value = 0.9562188673846281
$ws.Range("M7").Value = 121.928739
$ws.Range("N7").Value = 365.786217
$ws.Range("O7").Value = 0.2282232151508951
$ws.Range("P7").Value = 0.2419720431319445
$ws.Range("Q7").Value = 3654.009343776339
$ws.Range("R7").Value = 32886.08409398705
$ws.Range("S7").Value = 0.2143469700863104
$ws.Range("T7").Value = 0.2313782330223724
$ws.Range("I8").Value = 0.9391988012463586
$ws.Range("J8").Value = 0.9562188673846281
$ws.Range("O8").Value = 0.2768624053389947
$ws.Range("P8").Value = 0.2935413991166814
$ws.Range("R8").Value = 39894.80359577697
$ws.Range("S8").Value = 0.2600288392045673
$ws.Range("T8").Value = 0.2806898241938522
$ws.Range("I9").Value = 0.9391988012463586
$ws.Range("J9").Value = 0.9562188673846281
$ws.Range("M9").Value = 83.50496933333334
$ws.Range("N9").Value = 250.514908
$ws.Range("O9").Value = 0.1563025480180701
$ws.Range("P9").Value = 0.1657186665504434
$ws.Range("Q9").Value = 2502.510406474036
$ws.Range("R9").Value = 22522.59365826632
$ws.Range("S9").Value = 0.1467991657303228
$ws.Range("T9").Value = 0.1584633156333558
$ws.Range("I10").Value = 0.9391988012463586
$ws.Range("J10").Value = 0.9562188673846281
$ws.Range("M10").Value = 91.06846250000001
$ws.Range("N10").Value = 182.136925
$ws.Range("O10").Value = 0.1704597085236707
$ws.Range("P10").Value = 0.1204857969594293
$ws.Range("Q10").Value = 2729.176202653463
$ws.Range("R10").Value = 16375.05721592078
$ws.Range("S10").Value = 0.1600955539062353
$ws.Range("T10").Value = 0.1152107923044798
$ws.Range("I11").Value = 0.9391988012463586
$ws.Range("J11").Value = 0.9562188673846281
$ws.Range("M11").Value = 89.83562999999999
$ws.Range("N11").Value = 269.50689
$ws.Range("O11").Value = 0.1681521229683693
$ws.Range("P11").Value = 0.1782820942415013
$ws.Range("Q11").Value = 2692.23018392763
$ws.Range("R11").Value = 24230.07165534867
$ws.Range("S11").Value = 0.1579282723189227
$ws.Range("T11").Value = 0.1704767022305679
$ws.Range("I12").Value = 0.003433215419517052
$ws.Range("J12").Value = 0.003495431803768782
$ws.Range("M12").Value = 121.928739
$ws.Range("N12").Value = 365.786217
$ws.Range("O12").Value = 0.2282232151508951
$ws.Range("P12").Value = 0.2419720431319445
$ws.Range("Q12").Value = 13.357130785798
$ws.Range("R12").Value = 120.214177072182
$ws.Range("S12").Value = 0.0007835394613478107
$ws.Range("T12").Value = 0.0008457967751863103
$ws.Range("I13").Value = 0.003433215419517052
$ws.Range("J13").Value = 0.003495431803768782
$ws.Range("O13").Value = 0.2768624053389947
$ws.Range("P13").Value = 0.2935413991166814
$ws.Range("S13").Value = 0.0009505282790944169
$ws.Range("T13").Value = 0.001026053942195234
$ws.Range("I14").Value = 0.003433215419517052
$ws.Range("J14").Value = 0.003495431803768782
$ws.Range("M14").Value = 83.50496933333334
$ws.Range("N14").Value = 250.514908
$ws.Range("O14").Value = 0.1563025480180701
$ws.Range("P14").Value = 0.1657186665504434
$ws.Range("Q14").Value = 9.147858050507557
$ws.Range("R14").Value = 82.33072245456799
$ws.Range("S14").Value = 0.0005366203179654426
$ws.Range("T14").Value = 0.0005792582975385736
$ws.Range("I15").Value = 0.003433215419517052
$ws.Range("J15").Value = 0.003495431803768782
$ws.Range("M15").Value = 91.06846250000001
$ws.Range("N15").Value = 182.136925
$ws.Range("O15").Value = 0.1704597085236707
$ws.Range("P15").Value = 0.1204857969594293
$ws.Range("Q15").Value = 9.976428642258334
$ws.Range("R15").Value = 59.85857185355
$ws.Range("S15").Value = 0.0005852248997098487
$ws.Range("T15").Value = 0.0004211498865944173
$ws.Range("I16").Value = 0.003433215419517052
$ws.Range("J16").Value = 0.003495431803768782
$ws.Range("M16").Value = 89.83562999999999
$ws.Range("N16").Value = 269.50689
$ws.Range("O16").Value = 0.1681521229683693
$ws.Range("P16").Value = 0.1782820942415013
$ws.Range("Q16").Value = 9.84137348566
$ws.Range("R16").Value = 88.57236137094
$ws.Range("S16").Value = 0.000577302461399533
$ws.Range("T16").Value = 0.0006231729022542469
$ws.Range("G17").Value = 1.7038495
$ws.Range("H17").Value = 3.407699
$ws.Range("I17").Value = 0.05339802440257681
$ws.Range("J17").Value = 0.03624379868390631
$ws.Range("M17").Value = 121.928739
$ws.Range("N17").Value = 365.786217
$ws.Range("O17").Value = 0.2282232151508951
$ws.Range("P17").Value = 0.2419720431319445
$ws.Range("Q17").Value = 207.7482209807805
$ws.Range("R17").Value = 1246.489325884683
$ws.Range("S17").Value = 0.01218666881186203
$ws.Range("T17").Value = 0.008769986018407692
$ws.Range("G18").Value = 1.7038495
$ws.Range("H18").Value = 3.407699
$ws.Range("I18").Value = 0.05339802440257681
$ws.Range("J18").Value = 0.03624379868390631
$ws.Range("O18").Value = 0.2768624053389947
$ws.Range("P18").Value = 0.2935413991166814
$ws.Range("Q18").Value = 252.023757213335
$ws.Range("R18").Value = 1512.14254328001
$ws.Range("S18").Value = 0.01478390547644775
$ws.Range("T18").Value = 0.0106390553749772
$ws.Range("G19").Value = 1.7038495
$ws.Range("H19").Value = 3.407699
$ws.Range("I19").Value = 0.05339802440257681
$ws.Range("J19").Value = 0.03624379868390631
$ws.Range("M19").Value = 83.50496933333334
$ws.Range("N19").Value = 250.514908
$ws.Range("O19").Value = 0.1563025480180701
$ws.Range("P19").Value = 0.1657186665504434
$ws.Range("Q19").Value = 142.2799002461153
$ws.Range("R19").Value = 853.679401476692
$ws.Range("S19").Value = 0.008346247273253839
$ws.Range("T19").Value = 0.006006273988619669
$ws.Range("G20").Value = 1.7038495
$ws.Range("H20").Value = 3.407699
$ws.Range("I20").Value = 0.05339802440257681
$ws.Range("J20").Value = 0.03624379868390631
$ws.Range("M20").Value = 91.06846250000001
$ws.Range("N20").Value = 182.136925
$ws.Range("O20").Value = 0.1704597085236707
$ws.Range("P20").Value = 0.1204857969594293
$ws.Range("Q20").Value = 155.1669542963938
$ws.Range("R20").Value = 620.667817185575
$ws.Range("S20").Value = 0.0091022116754031
$ws.Range("T20").Value = 0.004366862969267568
$ws.Range("G21").Value = 1.7038495
$ws.Range("H21").Value = 3.407699
$ws.Range("I21").Value = 0.05339802440257681
$ws.Range("J21").Value = 0.03624379868390631
$ws.Range("M21").Value = 89.83562999999999
$ws.Range("N21").Value = 269.50689
$ws.Range("O21").Value = 0.1681521229683693
$ws.Range("P21").Value = 0.1782820942415013
$ws.Range("Q21").Value = 153.066393257685
$ws.Range("R21").Value = 918.39835954611
$ws.Range("S21").Value = 0.00897899116561008
$ws.Range("T21").Value = 0.006461620332634186
